$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 231: update existing date (2023-01-01 -> 2023-01-31)
$ws.Range("A231").Value = 44957

# Row 232: update existing date (2023-02-01 -> 2023-02-28)
$ws.Range("A232").Value = 44985

# Row 233: new monthly entry (2023-03-31), 1.25 earned
$ws.Range("A233").Value = 45016
$ws.Range("C233").Value = 1.25

# Row 234: new monthly entry (2023-04-30), 1.25 earned
$ws.Range("A234").Value = 45046
$ws.Range("C234").Value = 1.25

# Row 235: new monthly entry (2023-05-31), 1.25 earned
$ws.Range("A235").Value = 45077
$ws.Range("C235").Value = 1.25

# Row 236: new monthly entry (2023-06-30), 1.25 earned
$ws.Range("A236").Value = 45107
$ws.Range("C236").Value = 1.25

# Row 237: new monthly entry (2023-07-31), 1.25 earned
$ws.Range("A237").Value = 45138
$ws.Range("C237").Value = 1.25

# Row 238: new monthly entry (2023-08-31), 1.25 earned
$ws.Range("A238").Value = 45169
$ws.Range("C238").Value = 1.25

# Row 239: new monthly entry (2023-09-30), 1.25 earned
$ws.Range("A239").Value = 45199
$ws.Range("C239").Value = 1.25

# Row 240: new monthly entry (2023-10-31), with FL(3-0-00) particular, 1.25 earned,
# 3 days absence w/pay, remarks "10/23-25/2023"
$ws.Range("A240").Value = 45230
$ws.Range("B240").Value = "FL(3-0-00)"
$ws.Range("C240").Value = 1.25
$ws.Range("D240").Value = 3
$ws.Range("K240").Value = "10/23-25/2023"

# Row 241: new monthly entry (2023-11-30), with SP(1-0-00) particular, 1.25 earned,
# remarks date 12/4/2023
$ws.Range("A241").Value = 45260
$ws.Range("B241").Value = "SP(1-0-00)"
$ws.Range("C241").Value = 1.25
$ws.Range("K241").Value = 45264

# Row 242: new monthly entry (2023-12-31), with FL(2-0-0) particular,
# 2 days absence w/pay, remarks "12/5,6/2023"
$ws.Range("A242").Value = 45291
$ws.Range("B242").Value = "FL(2-0-0)"
$ws.Range("D242").Value = 2
$ws.Range("K242").Value = "12/5,6/2023"

# Row 243: new monthly entry (2024-01-31)
$ws.Range("A243").Value = 45322

# Row 244: new monthly entry (2024-02-29)
$ws.Range("A244").Value = 45351

# Row 245: new monthly entry (2024-03-31)
$ws.Range("A245").Value = 45382
